$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 154
$date = 43936

$data = @(
    @("Helse Midt-Norge", 4),
    @("Helse Nord", 5),
    @("Helse Sør-Øst", 35),
    @("Helse Vest", 8)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $name = $data[$i][0]
    $value = $data[$i][1]

    $ws.Cells.Item($row, 1).Value = $date
    $ws.Cells.Item($row, 1).NumberFormat = "yyyy-mm-dd"
    $ws.Cells.Item($row, 2).Value = $name
    $ws.Cells.Item($row, 3).Value = $value
}
